$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers: I1 ("I0") and J1 ("IF") ---
# Copy formatting from an existing header cell (H1 uses the bold/centered/
# bordered header style) so the new header cells pick up the same style
# index instead of minting a new one, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- New data columns I and J for rows 2-31 ---
$data = @(
    @(1, 2),
    @(4, 7),
    @(5, 7),
    @(6, 7),
    @(8, 8),
    @(9, 9),
    @(3, 5),
    @(4, 5),
    @(7, 7),
    @(5, 6),
    @(7, 8),
    @(7, 7),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(7, 8),
    @(10, 10),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(4, 5),
    @(8, 9),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(5, 5)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
